$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# All edited "D" (Price) cells get a temporary text number-format so Excel
# does not auto-coerce the numeric-looking string into a binary double
# (which would introduce float rounding noise like 312.57999999999998).
# The format is reset back to the default "Normal" style right after the
# value is written so the saved cell XML carries no stray style index.

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '43.900.26'
$ws.Range('D2').Style = "Normal"
$ws.Range('E2').Value = '  -0.09%  '

$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '2.233.60'
$ws.Range('D3').Style = "Normal"
$ws.Range('E3').Value = '  -1.73%  '

$ws.Range('E4').Value = '  +0.38%  '

$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '312.58'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  -2.70%  '

$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '98.34'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  -4.56%  '

$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.567'
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').Value = '  -3.58%  '

$ws.Range('E8').Value = '  +0.29%  '

$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.532'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  -7.09%  '

$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '36.02'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  -5.31%  '

$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.0818'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  -2.86%  '

$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '7.32'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  -6.91%  '

$ws.Range('E13').Value = '  -3.26%  '

$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '2.573.92'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  -1.63%  '

$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '2.241.05'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  -1.50%  '

$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '0.833'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  -5.16%  '

$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '14.05'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  -3.66%  '

$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '43.775.08'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  -0.26%  '

$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '12.90'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  -9.81%  '

$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '0.0₃0958'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  -3.72%  '

$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '6.32'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  -5.55%  '

$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '64.67'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  -2.34%  '

$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '2.98'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  -6.95%  '

$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '232.43'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  -2.92%  '

$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '2.02'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  -9.94%  '

$ws.Range('E26').Value = '  +0.09%  '

$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '10.14'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  -0.89%  '

$ws.Range('E28').Value = '  -1.86%  '

$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '36.62'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  -6.53%  '

$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '5.89'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  -9.65%  '

$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '157.27'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  -2.12%  '

$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '19.87'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  -3.35%  '

$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '0.0824'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  -7.04%  '

$ws.Range('E34').Value = '  -1.54%  '

$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '3.15'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  -7.59%  '

$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.109'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  +2.43%  '

$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '1.89'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  -6.77%  '

$ws.Range('E38').Value = '  -4.33%  '

$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '15.74'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  +0.72%  '

$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '3.58'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  -7.86%  '

$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '4.04'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  -10.58%  '

$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.0306'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  -7.07%  '

$ws.Range('E43').Value = '  +0.29%  '

$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '1.714.01'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  -5.91%  '

$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.192'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  -7.89%  '

$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '80.00'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  -7.30%  '

$ws.Range('B47').Value = 'Stacks'
$ws.Range('C47').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '1.67'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  -0.24%  '

$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '5.08'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  -6.32%  '

$ws.Range('B49').Value = 'ordi'
$ws.Range('C49').Value = 'https://coinranking.com/coin/j7-7vPrOi+ordi-ordi'
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '72.69'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  -5.37%  '

$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '101.09'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  -3.09%  '

$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '56.05'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  -6.40%  '
